# Updates cryptocurrency Price (D) and Volume(1h) (E) columns with latest
# scraped values, per the "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "57.326.77"
$cell.ClearFormats()
$ws.Range("E2").Value = "  +0.58%  "

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.357.03"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.83%  "

# Row 4
$ws.Range("E4").Value = "  -0.60%  "

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "520.58"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -0.33%  "

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "135.59"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +0.11%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  +0.32%  "

# Row 9
$ws.Range("E9").Value = "  -1.67%  "

# Row 10
$ws.Range("E10").Value = "  +4.65%  "

# Row 11
$ws.Range("E11").Value = "  -0.65%  "

# Row 12
$ws.Range("E12").Value = "  -1.13%  "

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "24.35"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +0.97%  "

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.779.77"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.83%  "

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "57.321.51"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.40%  "

# Row 16
$ws.Range("E16").Value = "  -1.03%  "

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.371.89"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +1.39%  "

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "10.59"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "329.04"
$cell.ClearFormats()
$ws.Range("E19").Value = "  +1.63%  "

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "4.23"
$cell.ClearFormats()
$ws.Range("E20").Value = "  -1.71%  "

# Row 21
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "61.28"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +0.21%  "

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "8.78"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +12.45%  "

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.165"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +3.32%  "

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "0.995"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.59%  "

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "1.35"
$cell.ClearFormats()
$ws.Range("E27").Value = "  +10.09%  "

# Row 28
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0742"
$cell.ClearFormats()
$ws.Range("E28").Value = "  -1.37%  "

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "167.58"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -2.74%  "

# Row 30
$ws.Range("E30").Value = "  -0.12%  "

# Row 31
$ws.Range("E31").Value = "  -1.28%  "

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "18.55"
$cell.ClearFormats()
$ws.Range("E32").Value = "  +0.39%  "

# Row 33
$ws.Range("E33").Value = "  +0.03%  "

# Row 34
$ws.Range("E34").Value = "  +1.98%  "

# Row 35
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.994"
$cell.ClearFormats()
$ws.Range("E35").Value = "  +0.38%  "

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.926"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -3.54%  "

# Row 37
$ws.Range("E37").Value = "  -0.83%  "

# Row 38
$ws.Range("E38").Value = "  +5.41%  "

# Row 39
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "38.83"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +3.15%  "

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "150.05"
$cell.ClearFormats()
$ws.Range("E40").Value = "  +6.99%  "

# Row 41
$ws.Range("E41").Value = "  +0.32%  "

# Row 42
$ws.Range("E42").Value = "  +0.67%  "

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "5.34"
$cell.ClearFormats()
$ws.Range("E43").Value = "  +2.36%  "

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "284.08"
$cell.ClearFormats()
$ws.Range("E44").Value = "  +1.74%  "

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0941"
$cell.ClearFormats()
$ws.Range("E45").Value = "  +0.90%  "

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0509"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -0.66%  "

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.565"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -0.71%  "

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "18.24"
$cell.ClearFormats()
$ws.Range("E48").Value = "  +4.30%  "

# Row 49
$ws.Range("E49").Value = "  +0.88%  "

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.386"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.61%  "

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "17.67"
$cell.ClearFormats()
$ws.Range("E51").Value = "  +3.39%  "
